$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.727.62"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "3.739.85"
$ws.Range("E3").Value = "  +18.92%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'613.94"
$ws.Range("E5").Value = "  +6.18%  "
$ws.Range("D6").Value = "'177.00"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").Value = "3.735.70"
$ws.Range("E7").Value = "  +18.84%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +10.34%  "
$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  +7.35%  "
$ws.Range("D13").Value = "'41.18"
$ws.Range("E13").Value = "  +11.42%  "
$ws.Range("D14").Value = "'0.0000257"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("D15").Value = "4.374.83"
$ws.Range("E15").Value = "  +19.37%  "
$ws.Range("D16").Value = "3.751.74"
$ws.Range("E16").Value = "  +19.41%  "
$ws.Range("D17").Value = "69.911.31"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "'7.63"
$ws.Range("E19").Value = "  +6.86%  "
$ws.Range("D20").Value = "'518.13"
$ws.Range("E20").Value = "  +6.04%  "
$ws.Range("D21").Value = "'16.77"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  +19.66%  "
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("D24").Value = "'88.85"
$ws.Range("E24").Value = "  +5.81%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("D26").Value = "'13.62"
$ws.Range("E26").Value = "  +4.47%  "
$ws.Range("D27").Value = "'10.97"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'0.0000127"
$ws.Range("E29").Value = "  +33.11%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  +6.01%  "
$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "'31.60"
$ws.Range("E33").Value = "  +11.70%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'6.25"
$ws.Range("E36").Value = "  +9.76%  "
$ws.Range("E37").Value = "  +7.84%  "
$ws.Range("D38").Value = "'0.342"
$ws.Range("E38").Value = "  +4.57%  "
$ws.Range("E39").Value = "  +6.56%  "
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").Value = "'51.47"
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("E42").Value = "  +5.47%  "
$ws.Range("D43").Value = "'44.47"
$ws.Range("E43").Value = "  -8.91%  "
$ws.Range("D44").Value = "'423.79"
$ws.Range("E44").Value = "  +5.15%  "
$ws.Range("D45").Value = "3.073.85"
$ws.Range("E45").Value = "  +9.53%  "
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'0.0367"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("D48").Value = "'27.92"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.52"
$ws.Range("E49").Value = "  +5.85%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'136.08"
$ws.Range("E50").Value = "  +0.70%  "
